# Apply commData.xlsx edits:
#  - Rename header E1 "AffectedPop" -> "VulPop"
#  - Remove the "MaxDistance" column (old column G, all "2000" values)
#    causing the old column H ("Remarks") to shift left into column G
#  - Dimension becomes A1:G17 as a consequence

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the AffectedPop header to VulPop
$ws.Range("E1").Value = "VulPop"

# Delete the MaxDistance column (column G); this shifts Remarks (H) left into G
$ws.Range("G1").EntireColumn.Delete()
